# Auto-applied update mirroring the "Updated cryptos list" GitHub Actions commit.
# Updates Price (D) and Volume(1h) (E) columns for many rows, and swaps the
# LidoDAOToken / Toncoin rows (25 and 26) to reflect the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column values that look numeric ("303.19", "1.000", ...) must be forced
# to Text format first, otherwise Excel auto-converts them to numbers (losing
# trailing zeros / introducing floating-point noise) when assigned via .Value.
$priceTextCells = @(
    "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14",
    "D15", "D18", "D19", "D21", "D22", "D23", "D25", "D26",
    "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36",
    "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45",
    "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $priceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values cell by cell (row order matches the sheet).
$ws.Range("D2").Value = "23.193.03"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.601.55"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "303.19"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "0.3783"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "52.12"
$ws.Range("E8").Value = "  +4.72%  "
$ws.Range("D9").Value = "0.3613"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "1.270"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "0.08112"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "22.69"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").Value = "6.576"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "7.409"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "1.603.16"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "94.04"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "0.06866"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "6.542"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "12.97"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "23.189.20"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.396"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.997"
$ws.Range("E26").Value = "  +10.57%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "149.38"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "5.245"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "133.92"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "2.382"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "6.767"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "1.780.34"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "0.9678"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "0.07496"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "10.29"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").Value = "0.02713"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").Value = "0.08799"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "6.092"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").Value = "0.7101"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "1.361"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").Value = "12.51"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "15.47"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").Value = "0.6524"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "2.312"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "4.020"
$ws.Range("D48").Value = "132.06"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "0.07966"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "1.199"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "1.210"
$ws.Range("E51").Value = "  +1.08%  "
